# Card8 sheet: column A (the "card" column) for the rows that were
# showing "2" should read "8" instead, matching the other rows in the
# table (e.g. rows 2 and 8 already contain "8"). The source data stores
# these as text, so force a Text number format before writing the value
# to keep "8" as a string rather than have it auto-convert to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card8")

$rows = @(3, 4, 5, 6, 7, 9, 10, 11, 12, 13)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = "8"
}
